$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date/time value for row 87 (column A)
$ws.Range("A87").Value = 45462.2916666667

# Add new row 88 with the new data point
$ws.Range("A88").Value = 45463.5694791667
$ws.Range("B88").Value = 1500
$ws.Range("C88").Value = 6
$ws.Range("D88").Value = 5.96000003814697
$ws.Range("E88").Value = 6
$ws.Range("F88").Value = 5.96000003814697
$ws.Range("G88").Formula = "=""5.96000003814697"""
$ws.Range("G88").Copy()
$ws.Range("G88").PasteSpecial(-4163)
$ws.Range("H88").Value = "PAL.MI"

# Match the date format style used in column A (style index 1 -> custom date format)
$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122)
